$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("A3SS")

$ws.Range("A2").Value = "chr10"
$ws.Range("B2").Value = "NCOA7"
$ws.Range("C2").Value = 30691755
$ws.Range("D2").Value = 30691787
$ws.Range("E2").Value = 0.8698350074869731
$ws.Range("F2").Value = 0.026160299348494
$ws.Range("G2").Value = 0.026160299348494
$ws.Range("H2").Value = 88.33333333333333
$ws.Range("I2").Value = 95.53465184700885
$ws.Range("J2").Value = 102.7272727272727
$ws.Range("K2").Value = 123.4163125204315
$ws.Range("L2").Value = "chr10:NCOA7:30691755-30691787"
